$d = $word.ActiveDocument

# 1. Change the title text
[void]$d.Content.Find.Execute("Questions I have about this Blog", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Things to do with Website", 2)

# 2. Strike through the "multiple pictures" bullet (2nd paragraph)
$d.Paragraphs(2).Range.Font.StrikeThrough = 1

# 3. Strike through the "div and state property for videos" bullet (4th paragraph)
$d.Paragraphs(4).Range.Font.StrikeThrough = 1

# 4. Append new bullet items after the last paragraph, using raw OOXML so we
#    can precisely control run-splitting / proofErr markers exactly as authored.
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$endRange = $d.Range($d.Content.End, $d.Content.End)

$newXml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Figure out how to fix right margin for Add-to-Cart box for </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>ProductScreen</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Fix Read More Button on The </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>BlogHome</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> Screen</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Change Colors to my liking.</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Add Social Media Icon Buttons</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Fix Styling on Orders Screen</w:t>
  </w:r>
</w:p>
"@

[void]$endRange.InsertXML($newXml)
